# Fill in the two data points on the "Relatório de Atividades" cover
# table that were left blank in the template:
#   1. The bank-account number gets its check-digit suffix "-7".
#   2. The "Vigência das Atividades do Projeto" line gets the actual
#      validity period of the project.

$d = $word.ActiveDocument

# 1. Conta Bancária: 10.738  ->  Conta Bancária: 10.738-7
$d.Content.Find.Execute("Conta Bancária: 10.738", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Conta Bancária: 10.738-7", 2) | Out-Null

# 2. Vigência das Atividades do Projeto:  ->  ... 23/02/2021 a 23/02/2022
$d.Content.Find.Execute("Vigência das Atividades do Projeto:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Vigência das Atividades do Projeto: 23/02/2021 a 23/02/2022", 2) | Out-Null
